# qs-pmfspdfscdfs.docx edit:
#  1. Collapse the word-by-word runs in the Title / Author / Abstract
#     paragraphs into single runs (text content is unchanged).
#  2. Clear the explicit "left" paragraph alignment on the compact
#     table-cell paragraphs so it falls back to the style default
#     (removes the redundant <w:jc w:val="left"/>).

$d = $word.ActiveDocument

# --- 1. Merge split runs back into a single run per paragraph ---------

function Merge-ParagraphRuns($paraIndex, $text) {
    $rng = $d.Paragraphs.Item($paraIndex).Range
    $rng.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, $text, 2) | Out-Null
}

Merge-ParagraphRuns 1 "Questions: PMFs, PDFs, and CDFs"
Merge-ParagraphRuns 2 "Sophie Chowgule"
Merge-ParagraphRuns 4 "A selection of questions to test your understanding of Probability Mass Functions (PMFs), Probability Density Functions (PDFs), and Cumulative Distribution Functions (CDFs)."

# --- 2. Drop the redundant explicit left-alignment on Compact paragraphs

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Compact") {
        $p.Alignment = 0
    }
}
